# Updated Todo with current project and added rough materials and roundtracker script
$wb = $excel.ActiveWorkbook
$wsProgramming = $wb.Worksheets.Item("Programming")
$wsArt = $wb.Worksheets.Item("Art")

# --- Content changes on the "Programming" sheet ---
# Floor Health row: progress capitalization fixed ("In progress" -> "In Progress")
$wsProgramming.Range("C8").Value = "In Progress"

# Round Tracking row: tackled by Zach, currently "In Progress" (roundtracker script)
$wsProgramming.Range("B13").Value = "Zach"
$wsProgramming.Range("C13").Value = "In Progress"

# --- Active sheet / selection changes ---
# Make "Programming" the active/selected sheet (was "Art")
$wsProgramming.Select()
$wsProgramming.Range("C14").Select()

# Update the no-longer-active "Art" sheet's stored selection
$wsArt.Range("A24").Select()

# Re-select "Programming" so it ends up as the active sheet/tab in the saved file
$wsProgramming.Select()
$wsProgramming.Range("C14").Select()
